$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7287194209349384
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 2938.103010863317
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("G2").Value = 3187.470217437252

$ws.Range("B3").Value = 0.3464964993005633
$ws.Range("C3").Value = 0.3375848360084654
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 246.9852506941017
$ws.Range("G3").Value = 248.3820648804257
